$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 636, shifting existing rows 636.. down by one
# (this grows the used range from A1:D677 to A1:D678, matching the new
# <dimension> and re-numbering every subsequent data row by +1).
$ws.Rows.Item(636).Insert()

# Column A holds plain text dates (e.g. "2026/01/12"), not real Excel date
# serials. Force text formatting before the assignment so Excel's input
# auto-detection doesn't silently convert the string into a date value;
# then restore the default "Normal" style so the cell doesn't end up with
# a stray explicit format like its neighbours.
$ws.Range("A636").NumberFormat = "@"
$ws.Range("A636").Value = "2026/01/12"
$ws.Range("A636").Style = "Normal"

$ws.Range("B636").Value = "月"
$ws.Range("C636").Value = 19
$ws.Range("D636").Value = 181
